$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated accuracy values (data preprocessing change)
$ws.Range("C6").Value = 67.702241224180995
$ws.Range("C7").Value = 90.7938214128976
$ws.Range("C8").Value = -24.3945651726096
$ws.Range("C9").Value = 77.838711222724399
$ws.Range("C12").Value = 88.673858191680395

# Update the active selection to match the saved cursor position
$ws.Range("C15").Select()
